$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 830.06665

$ws.Range("H40").Value = 1218.1818
$ws.Range("J40").Value = 1342.8572
$ws.Range("L40").Value = 1342.8572
$ws.Range("N40").Value = -1692.8572

$ws.Range("H43").Value = 4430.1
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 4787.625
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 4787.625
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -4925.625

$ws.Range("H58").Value = 2407.5
$ws.Range("J58").Value = 2603.3333
$ws.Range("L58").Value = 7809.999899999999
$ws.Range("N58").Value = -8109.999899999999

$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = -3752

$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = -3142

$ws.Range("H113").Value = 4735.3125
$ws.Range("I113").Value = 4305.4165
$ws.Range("J113").Value = 6025
$ws.Range("K113").Value = 4305.4165
$ws.Range("L113").Value = 6025
$ws.Range("M113").Value = -1051.4165
$ws.Range("N113").Value = -12533

$ws.Range("H115").Value = 1989.5454
$ws.Range("I115").Value = 442.5
$ws.Range("J115").Value = 2333.3333
$ws.Range("K115").Value = 1327.5
$ws.Range("L115").Value = 6999.999899999999
$ws.Range("M115").Value = 239.5
$ws.Range("N115").Value = -10133.9999

$ws.Range("H116").Value = 3601.0938
$ws.Range("I116").Value = 3010.45
$ws.Range("K116").Value = 3010.45
$ws.Range("M116").Value = 431.5500000000002

$ws.Range("H132").Value = 4350220.5
$ws.Range("I132").Value = 4880398.5
$ws.Range("K132").Value = 14641195.5
$ws.Range("M132").Value = -14638665.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 19021.572
$ws.Range("J37").Value = 21025.166
$ws.Range("L37").Value = 21025.166
$ws.Range("N37").Value = -21571.166

$ws.Range("H61").Value = 2182.4614
$ws.Range("I61").Value = 1240.0667
$ws.Range("J61").Value = 5323.778
$ws.Range("K61").Value = 1240.0667
$ws.Range("L61").Value = 5323.778
$ws.Range("M61").Value = -1028.0667
$ws.Range("N61").Value = -5747.778

$ws.Range("H132").Value = 16395813
$ws.Range("I132").Value = 25642356
$ws.Range("K132").Value = 76927068
$ws.Range("M132").Value = -76924538

$ws.Range("H136").Value = 2182.4614
$ws.Range("I136").Value = 1240.0667
$ws.Range("J136").Value = 5323.778
$ws.Range("K136").Value = 3720.2001
$ws.Range("L136").Value = 15971.334
$ws.Range("M136").Value = -1170.2001
$ws.Range("N136").Value = -21071.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2440.5
$ws.Range("I99").Value = 1462.3334
$ws.Range("J99").Value = 5375
$ws.Range("K99").Value = 1462.3334
$ws.Range("L99").Value = 5375
$ws.Range("M99").Value = 35.66660000000002
$ws.Range("N99").Value = -8371

$ws.Range("H134").Value = 2311.7917
$ws.Range("I134").Value = 1472.0526
$ws.Range("K134").Value = 4416.1578
$ws.Range("M134").Value = -1881.1578

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4121.4443
$ws.Range("I99").Value = 1956
$ws.Range("K99").Value = 1956
$ws.Range("M99").Value = -458

$ws.Range("H126").Value = 4121.4443
$ws.Range("I126").Value = 1956
$ws.Range("K126").Value = 5868
$ws.Range("M126").Value = -3398

$ws.Range("H132").Value = 2053.532
$ws.Range("I132").Value = 1532
$ws.Range("K132").Value = 4596
$ws.Range("M132").Value = -2066

$ws.Range("H134").Value = 1913.3214
$ws.Range("I134").Value = 988.2
$ws.Range("K134").Value = 2964.6
$ws.Range("M134").Value = -429.6000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 8392.857
$ws.Range("I34").Value = 266.66666
$ws.Range("J34").Value = 14487.5
$ws.Range("K34").Value = 799.9999799999999
$ws.Range("L34").Value = 43462.5
$ws.Range("M34").Value = -715.9999799999999
$ws.Range("N34").Value = -43630.5

$ws.Range("H55").Value = 2053.3333
$ws.Range("I55").Value = 460
$ws.Range("J55").Value = 4045
$ws.Range("K55").Value = 1380
$ws.Range("L55").Value = 12135
$ws.Range("M55").Value = -1203
$ws.Range("N55").Value = -12489

$ws.Range("H87").Value = 10412.333
$ws.Range("I87").Value = 5008
$ws.Range("J87").Value = 15816.667
$ws.Range("K87").Value = 15024
$ws.Range("L87").Value = 47450.001
$ws.Range("M87").Value = -13776
$ws.Range("N87").Value = -49946.001

$ws.Range("H90").Value = 10412.333
$ws.Range("I90").Value = 5008
$ws.Range("J90").Value = 15816.667
$ws.Range("K90").Value = 45072
$ws.Range("L90").Value = 142350.003
$ws.Range("M90").Value = -38832
$ws.Range("N90").Value = -154830.003

$ws.Range("H113").Value = 1299767.6
$ws.Range("I113").Value = 9091159
$ws.Range("J113").Value = 1202.4166
$ws.Range("K113").Value = 27273477
$ws.Range("L113").Value = 3607.2498
$ws.Range("M113").Value = -27271307
$ws.Range("N113").Value = -7947.2498

$ws.Range("H120").Value = 19117.428
$ws.Range("I120").Value = 17940.666
$ws.Range("K120").Value = 53821.99800000001
$ws.Range("M120").Value = -48983.99800000001

$ws.Range("H123").Value = 1578.4
$ws.Range("I123").Value = 973
$ws.Range("K123").Value = 2919
$ws.Range("M123").Value = -469

$ws.Range("H139").Value = 8622934
$ws.Range("I139").Value = 13890689
$ws.Range("J139").Value = 2969.9092
$ws.Range("K139").Value = 41672067
$ws.Range("L139").Value = 8909.7276
$ws.Range("M139").Value = -41666927
$ws.Range("N139").Value = -19189.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10499.875
$ws.Range("I5").Value = 6333.3335
$ws.Range("K5").Value = 6333.3335
$ws.Range("M5").Value = -6221.3335

$ws.Range("H22").Value = 47577.715
$ws.Range("I22").Value = 20004
$ws.Range("K22").Value = 20004
$ws.Range("M22").Value = -19475

$ws.Range("H132").Value = 3171.4187
$ws.Range("I132").Value = 2910.92
$ws.Range("K132").Value = 8732.76
$ws.Range("M132").Value = -6202.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2777.7778
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 3866.6667
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 3866.6667
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -4242.6667

$ws.Range("H55").Value = 960.8125
$ws.Range("I55").Value = 225.125
$ws.Range("J55").Value = 1696.5
$ws.Range("K55").Value = 225.125
$ws.Range("L55").Value = 1696.5
$ws.Range("M55").Value = -52.125
$ws.Range("N55").Value = -2042.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23333.334
$ws.Range("J123").Value = 23333.334
$ws.Range("L123").Value = 23333.334
$ws.Range("N123").Value = -33133.334

Write-Host "Updated all Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with refreshed market data."
